$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Updated dSF (column F) values, repulled from source data.
$values = @{
    2  = -2
    3  = 2
    4  = -1
    6  = 2
    7  = -1
    8  = 1
    9  = -3
    11 = 3
    12 = 4
    13 = -4
    14 = -2
    15 = -3
    16 = 8
    17 = 4
    18 = -2
    19 = 5
    20 = 6
    21 = 1
    22 = 3
    23 = 1
    24 = -6
    25 = 3
    26 = 3
    27 = 2
    28 = -1
    29 = 12
    30 = -3
    32 = 2
    33 = 1
    34 = 1
    35 = 2
    36 = 2
}

foreach ($row in $values.Keys) {
    $ws.Range("F$row").Value = $values[$row]
}
